# Update the "想去人数" (attendance/interest count) figures for two events
# across the "展览" and "全部类型" worksheets, matching the upstream data
# refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3426
$ws1.Range("F5").Value = 6977
$ws1.Range("F6").Value = 2459
$ws1.Range("F7").Value = 47
$ws1.Range("F14").Value = 574

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3426
$ws4.Range("F6").Value = 6978
$ws4.Range("F7").Value = 2459
$ws4.Range("F8").Value = 47
$ws4.Range("F15").Value = 574
